$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labels = @(
    "Agri. & Vet sciences",
    "Built Env. & Design",
    "Technology",
    "Earth sciences",
    "Env. sciences",
    "Physical sciences",
    "Philo & Relig",
    "Law & Legal",
    "Chemical sciences",
    "Comm. Manage. Tourism",
    "Education",
    "Creat. Arts & Writing",
    "Math sciences",
    "Pysch. & Cognit. sciences",
    "Biological sciences",
    "Economics",
    "Inf. & Comp. sciences",
    "Engineering",
    "Hist. & Archaeology",
    "Lang. Comms. & Culture",
    "Studies Human Society",
    "Med. & Health Sciences"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $labels[$i]
}
